$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.8700167400677685
$ws.Range("B3").Value = 0.8557983477857647
$ws.Range("B4").Value = 0.8477439638064084
$ws.Range("B5").Value = 0.8446317010680957
$ws.Range("B6").Value = 0.8441251789210469
$ws.Range("B7").Value = 0.8477013025128031
$ws.Range("B8").Value = 0.8649739950075457
$ws.Range("B9").Value = 0.904209854685206
$ws.Range("B10").Value = 0.9363154015657358
$ws.Range("B11").Value = 0.9516354510921019
$ws.Range("B12").Value = 0.9575396873666762
$ws.Range("B13").Value = 0.9562635287275327
$ws.Range("B14").Value = 0.9521191343189344
$ws.Range("B15").Value = 0.9495939705872729
$ws.Range("B16").Value = 0.9353285925649857
$ws.Range("B17").Value = 0.9267604179175066
$ws.Range("B18").Value = 0.921899533835898
$ws.Range("B19").Value = 0.9202652803807609
$ws.Range("B20").Value = 0.9276655488540655
$ws.Range("B21").Value = 0.9533336513056838
$ws.Range("B22").Value = 0.9707088467760627
$ws.Range("B23").Value = 0.9613804985338277
$ws.Range("B24").Value = 0.9272561364640239
$ws.Range("B25").Value = 0.8930203668474235

$ws.Range("C2").Value = 0.06848206782430566
$ws.Range("C3").Value = 0.0612012889496043
$ws.Range("C4").Value = 0.05674671088911509
$ws.Range("C5").Value = 0.05493530132503111
$ws.Range("C6").Value = 0.05463474827404013
$ws.Range("C7").Value = 0.0567222660832698
$ws.Range("C8").Value = 0.06596829213167155
$ws.Range("C9").Value = 0.08423172162511605
$ws.Range("C10").Value = 0.09774027504062133
$ws.Range("C11").Value = 0.103907526197645
$ws.Range("C12").Value = 0.1062462406610791
$ws.Range("C13").Value = 0.1057424074083144
$ws.Range("C14").Value = 0.1040998665731649
$ws.Range("C15").Value = 0.1030941977487316
$ws.Range("C16").Value = 0.09733768700280621
$ws.Range("C17").Value = 0.0938120148095436
$ws.Range("C18").Value = 0.09178621795379627
$ws.Range("C19").Value = 0.09110067076548489
$ws.Range("C20").Value = 0.09418711280662251
$ws.Range("C21").Value = 0.1045822297808741
$ws.Range("C22").Value = 0.1113953939604073
$ws.Range("C23").Value = 0.1077572670931772
$ws.Range("C24").Value = 0.09401752740274105
$ws.Range("C25").Value = 0.07927582465372041

$ws.Range("E2").Value = 0.04685825844356817
$ws.Range("E3").Value = 0.04698995929574146
$ws.Range("E4").Value = 0.04708437649157204
$ws.Range("E5").Value = 0.04712626947171117
$ws.Range("E6").Value = 0.0471334324364534
$ws.Range("E7").Value = 0.04708492762628325
$ws.Range("E8").Value = 0.04690086166632312
$ws.Range("E9").Value = 0.0466470478207448
$ws.Range("E10").Value = 0.04652538467685652
$ws.Range("E11").Value = 0.0464840077758808
$ws.Range("E12").Value = 0.04647033934260136
$ws.Range("E13").Value = 0.0464731942639478
$ws.Range("E14").Value = 0.0464828432290787
$ws.Range("E15").Value = 0.04648901371988856
$ws.Range("E16").Value = 0.04652836910165359
$ws.Range("E17").Value = 0.04655608387174848
$ws.Range("E18").Value = 0.04657333993796087
$ws.Range("E19").Value = 0.04657940870452837
$ws.Range("E20").Value = 0.04655299752431041
$ws.Range("E21").Value = 0.04647995487945522
$ws.Range("E22").Value = 0.04644387152208651
$ws.Range("E23").Value = 0.04646206631586125
$ws.Range("E24").Value = 0.04655438874204698
$ws.Range("E25").Value = 0.0467042933585553

$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("F25").Value = 0.5279251897347166

$ws.Range("G2").Value = 0.002630660161191975
$ws.Range("G3").Value = 0.002634995088196433
$ws.Range("G4").Value = 0.002637795453142711
$ws.Range("G5").Value = 0.0026389716177181
$ws.Range("G6").Value = 0.002639169035883274
$ws.Range("G7").Value = 0.00263781117367604
$ws.Range("G8").Value = 0.002632126125791998
$ws.Range("G9").Value = 0.002622073043274365
$ws.Range("G10").Value = 0.002615347367184607
$ws.Range("G11").Value = 0.002612429489765896
$ws.Range("G12").Value = 0.002611344817012722
$ws.Range("G13").Value = 0.002611577521161933
$ws.Range("G14").Value = 0.002612339847533622
$ws.Range("G15").Value = 0.002612809430337407
$ws.Range("G16").Value = 0.002615540898995444
$ws.Range("G17").Value = 0.002617252774066436
$ws.Range("G18").Value = 0.002618250740264605
$ws.Range("G19").Value = 0.002618590928799236
$ws.Range("G20").Value = 0.002617069162241364
$ws.Range("G21").Value = 0.002612115384718189
$ws.Range("G22").Value = 0.002608995870869082
$ws.Range("G23").Value = 0.002610650046438318
$ws.Range("G24").Value = 0.002617152130256875
$ws.Range("G25").Value = 0.002624676177298065

$ws.Range("I2").Value = 3.231981165218542
$ws.Range("I3").Value = 3.156941376114759
$ws.Range("I4").Value = 3.111109536730723
$ws.Range("I5").Value = 3.092492504528138
$ws.Range("I6").Value = 3.089404736208934
$ws.Range("I7").Value = 3.110858220426877
$ws.Range("I8").Value = 3.206056252023316
$ws.Range("I9").Value = 3.394732407140609
$ws.Range("I10").Value = 3.534671743590565
$ws.Range("I11").Value = 3.598643214586247
$ws.Range("I12").Value = 3.622914023294442
$ws.Range("I13").Value = 3.617684795150296
$ws.Range("I14").Value = 3.600639055343379
$ws.Range("I15").Value = 3.590204102323753
$ws.Range("I16").Value = 3.530497452826751
$ws.Range("I17").Value = 3.493950274964845
$ws.Range("I18").Value = 3.472958629553176
$ws.Range("I19").Value = 3.465856216512336
$ws.Range("I20").Value = 3.49783774016538
$ws.Range("I21").Value = 3.605644538215614
$ws.Range("I22").Value = 3.676372402397504
$ws.Range("I23").Value = 3.638598491581575
$ws.Range("I24").Value = 3.49608015551965
$ws.Range("I25").Value = 3.343466383575318

$ws.Range("K2").Value = 0.7815826091351994
$ws.Range("K3").Value = 0.7580917126121562
$ws.Range("K4").Value = 0.7442437698174444
$ws.Range("K5").Value = 0.7387451996900722
$ws.Range("K6").Value = 0.7378408943189356
$ws.Range("K7").Value = 0.7441690289786891
$ws.Range("K8").Value = 0.7733633867766514
$ws.Range("K9").Value = 0.8351917697416695
$ws.Range("K10").Value = 0.88343087482545
$ws.Range("K11").Value = 0.9059925351012907
$ws.Range("K12").Value = 0.9146251605492068
$ws.Range("K13").Value = 0.9127620075641119
$ws.Range("K14").Value = 0.9067009610843115
$ws.Range("K15").Value = 0.9029999942701465
$ws.Range("K16").Value = 0.8819688610173273
$ws.Range("K17").Value = 0.8692252432562384
$ws.Range("K18").Value = 0.8619535545504107
$ws.Range("K19").Value = 0.8595014582455747
$ws.Range("K20").Value = 0.8705758080660644
$ws.Range("K21").Value = 0.9084788196857119
$ws.Range("K22").Value = 0.9337695843495908
$ws.Range("K23").Value = 0.9202238656777695
$ws.Range("K24").Value = 0.8699650470798019
$ws.Range("K25").Value = 0.8179730980558872

$ws.Range("L2").Value = 0.2456611017818062
$ws.Range("L3").Value = 0.2423689574970993
$ws.Range("L4").Value = 0.240455488757199
$ws.Range("L5").Value = 0.2397028794879503
$ws.Range("L6").Value = 0.2395795490060166
$ws.Range("L7").Value = 0.2404452288863155
$ws.Range("L8").Value = 0.2445035738363686
$ws.Range("L9").Value = 0.2533186942374215
$ws.Range("L10").Value = 0.2603190731884126
$ws.Range("L11").Value = 0.2636179223774917
$ws.Range("L12").Value = 0.2648835686769502
$ws.Range("L13").Value = 0.2646102577281368
$ws.Range("L14").Value = 0.2637217182913787
$ws.Range("L15").Value = 0.2631796037487959
$ws.Range("L16").Value = 0.2601057865891363
$ws.Range("L17").Value = 0.2582493857279644
$ws.Range("L18").Value = 0.2571923952827859
$ws.Range("L19").Value = 0.256836364873223
$ws.Range("L20").Value = 0.2584458889109698
$ws.Range("L21").Value = 0.2639822575681023
$ws.Range("L22").Value = 0.2676964424620252
$ws.Range("L23").Value = 0.2657053412261945
$ws.Range("L24").Value = 0.2583570178051673
$ws.Range("L25").Value = 0.2508420880240863

$ws.Range("M2").Value = 0.2391894589264361
$ws.Range("M3").Value = 0.2356253238422568
$ws.Range("M4").Value = 0.233591186121437
$ws.Range("M5").Value = 0.2328010626961046
$ws.Range("M6").Value = 0.2326722074765151
$ws.Range("M7").Value = 0.2335803731026722
$ws.Range("M8").Value = 0.2379285362177974
$ws.Range("M9").Value = 0.2476794584455106
$ws.Range("M10").Value = 0.2555914369610122
$ws.Range("M11").Value = 0.2593536934601914
$ws.Range("M12").Value = 0.260801825869045
$ws.Range("M13").Value = 0.2604889017888325
$ws.Range("M14").Value = 0.25947236225484
$ws.Range("M15").Value = 0.2588527560309259
$ws.Range("M16").Value = 0.2553488457480881
$ws.Range("M17").Value = 0.2532410726993106
$ws.Range("M18").Value = 0.2520440850452488
$ws.Range("M19").Value = 0.2516414415967176
$ws.Range("M20").Value = 0.2534638603226469
$ws.Range("M21").Value = 0.2597703081886422
$ws.Range("M22").Value = 0.2640286079132039
$ws.Range("M23").Value = 0.2617433672595197
$ws.Range("M24").Value = 0.2533630919410257
$ws.Range("M25").Value = 0.2449103513418081
